$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 28
$ws.Cells.Item($row, 1).Value = 1111142
$ws.Cells.Item($row, 2).Value = "RAJ SHUKLA"
$ws.Cells.Item($row, 3).Value = 30
$ws.Cells.Item($row, 4).Value = "MALE"
$ws.Cells.Item($row, 5).Value = 22093852093
$ws.Cells.Item($row, 6).Value = "professor.el175@gmail.com"
$ws.Cells.Item($row, 7).Value = "RT-PCR"
$ws.Cells.Item($row, 8).Value = "DINESH SHAH"
$ws.Cells.Item($row, 9).Value = 500
$ws.Cells.Item($row, 10).Value = "22/03/2022"
$ws.Cells.Item($row, 11).Value = "09:03 PM"
